# Updated user stories and functional reqs. for user pages
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the dependency value for REQ_ID 26 (row 30): "?" -> "None"
$ws.Range("D30").Value = "None"

# Append new functional requirements / user stories for list management on the user page
$newRows = @(
    @{ Row = 67; A = "19"; B = "63"; C = "User Page"; D = "47"; E = "Ability to delete a movie list from the movie list view tab" },
    @{ Row = 68; A = "19"; B = "64"; C = "User Page"; D = "48"; E = "Ability to delete a people list from the people list view tab" },
    @{ Row = 69; A = "19"; B = "65"; C = "User Page"; D = "47"; E = "Ability to delete a particular movie from a movie list" },
    @{ Row = 70; A = "19"; B = "66"; C = "User Page"; D = "48"; E = "Ability to delete a particular person from a person list" },
    @{ Row = 71; A = "19"; B = "67"; C = "User Page"; D = "47"; E = "Ability to add a movie to a movie list from within the movie list tab" },
    @{ Row = 72; A = "19"; B = "68"; C = "User Page"; D = "48"; E = "Abilility to add a person to a person list from within the person list tab" }
)

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("E" + $r.Row).NumberFormat = "@"
}

$ws.Range("E72").Select()
